# Update "想去人数" (column F) values on 展览, 演出 and 全部类型 sheets
# to reflect a refreshed data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 130
$ws1.Cells.Item(4, 6).Value = 881
$ws1.Cells.Item(9, 6).Value = 11912
$ws1.Cells.Item(11, 6).Value = 2138
$ws1.Cells.Item(16, 6).Value = 1218
$ws1.Cells.Item(17, 6).Value = 186
$ws1.Cells.Item(19, 6).Value = 752
$ws1.Cells.Item(22, 6).Value = 2914
$ws1.Cells.Item(23, 6).Value = 747
$ws1.Cells.Item(24, 6).Value = 3818
$ws1.Cells.Item(25, 6).Value = 3818
$ws1.Cells.Item(27, 6).Value = 834
$ws1.Cells.Item(31, 6).Value = 1019
$ws1.Cells.Item(32, 6).Value = 45
$ws1.Cells.Item(33, 6).Value = 90
$ws1.Cells.Item(38, 6).Value = 4378
$ws1.Cells.Item(40, 6).Value = 4494
$ws1.Cells.Item(41, 6).Value = 5514
$ws1.Cells.Item(43, 6).Value = 122
$ws1.Cells.Item(45, 6).Value = 169
$ws1.Cells.Item(46, 6).Value = 290
$ws1.Cells.Item(47, 6).Value = 76
$ws1.Cells.Item(49, 6).Value = 4103
$ws1.Cells.Item(50, 6).Value = 119

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(12, 6).Value = 859

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 130
$ws4.Cells.Item(6, 6).Value = 881
$ws4.Cells.Item(11, 6).Value = 11912
$ws4.Cells.Item(12, 6).Value = 2138
$ws4.Cells.Item(15, 6).Value = 1218
$ws4.Cells.Item(16, 6).Value = 186
$ws4.Cells.Item(19, 6).Value = 752
$ws4.Cells.Item(21, 6).Value = 747
$ws4.Cells.Item(22, 6).Value = 3818
$ws4.Cells.Item(26, 6).Value = 834
$ws4.Cells.Item(29, 6).Value = 1019
$ws4.Cells.Item(30, 6).Value = 45
$ws4.Cells.Item(31, 6).Value = 90
$ws4.Cells.Item(35, 6).Value = 4494
$ws4.Cells.Item(37, 6).Value = 122
$ws4.Cells.Item(38, 6).Value = 169
$ws4.Cells.Item(39, 6).Value = 290
$ws4.Cells.Item(43, 6).Value = 76
$ws4.Cells.Item(45, 6).Value = 4103
$ws4.Cells.Item(50, 6).Value = 119
